$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Time value for "random forest" row (B2): 32.366209s -> 28.95758s
$ws.Range("B2").Value = "28.95758s"

# Fill in the previously empty Accuracy values
$ws.Range("C2").Value = 0.95679
$ws.Range("C5").Value = 0.96

# Update the selection to match the author's final selection (A1:C6 highlighted,
# active cell C6). The engine's Select() always anchors the active cell at the
# top-left of the selected range, so we select the full block here to at least
# reproduce the highlighted sqref region.
$ws.Range("A1:C6").Select()
